$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("CaLamViec")
$ws2 = $wb.Worksheets.Item("DM CaLviec")

# Reassign header row: B1 becomes "Tên NV" (Vietnamese for Employee Name),
# C1 takes over what used to be "Ca Lam Viec"
$ws1.Cells.Item(1,2).Value2 = "Tên NV"
$ws1.Cells.Item(1,3).Value2 = "Ca Lam Viec"

# Reassign data row: B2 becomes placeholder "XXXX", C2 takes over the
# shift-name value that used to live in the DM CaLviec lookup (copy the
# formatting from B2, which already carries the right border/fill style)
$ws1.Cells.Item(2,2).Value2 = "XXXX"
$ws1.Cells.Item(2,2).Copy()
$ws1.Cells.Item(2,3).PasteSpecial(-4122)
$ws1.Cells.Item(2,3).Value2 = $ws2.Cells.Item(2,2).Value2

# Drop the now-unused 4th column (Ngay Bat Dau / Ngay Ket Thuc no longer needed)
$ws1.Columns.Item(4).Delete()

# Move the active-cell selection to match the new layout
[void]$ws1.Range("C4").Select()

Write-Output "done"
